$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps caused by re-sorting after the daily data refresh ---
$ws.Cells.Item(52, 1).Value = "Nigeria"
$ws.Cells.Item(53, 1).Value = "Armenia"
$ws.Cells.Item(70, 1).Value = "Sudan"
$ws.Cells.Item(71, 1).Value = "Noruega"
$ws.Cells.Item(85, 1).Value = "Gabon"
$ws.Cells.Item(86, 1).Value = "Etiopia"
$ws.Cells.Item(89, 1).Value = "Bulgaria"
$ws.Cells.Item(90, 1).Value = "Hungria"
$ws.Cells.Item(91, 1).Value = "Venezuela"
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"
$ws.Cells.Item(208, 1).Value = "Islas Malvinas"
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(212, 1).Value = "Seychelles"

# --- Updated case/death/recovered counts for the new data snapshot ---
$ws.Cells.Item(4, 2).Value = 2422671
$ws.Cells.Item(4, 3).Value = 34518
$ws.Cells.Item(4, 4).Value = 1017649
$ws.Cells.Item(4, 5).Value = 1281559
$ws.Cells.Item(4, 7).Value = 853
$ws.Cells.Item(4, 8).Value = 123463
$ws.Cells.Item(5, 4).Value = 613345
$ws.Cells.Item(5, 5).Value = 479912
$ws.Cells.Item(14, 2).Value = 192778
$ws.Cells.Item(14, 3).Value = 659
$ws.Cells.Item(14, 5).Value = 8092
$ws.Cells.Item(25, 2).Value = 73572
$ws.Cells.Item(25, 3).Value = 2389
$ws.Cells.Item(25, 4).Value = 30459
$ws.Cells.Item(25, 5).Value = 40709
$ws.Cells.Item(25, 7).Value = 94
$ws.Cells.Item(25, 8).Value = 2404
$ws.Cells.Item(52, 2).Value = 21371
$ws.Cells.Item(52, 3).Value = 452
$ws.Cells.Item(52, 4).Value = 7338
$ws.Cells.Item(52, 5).Value = 13500
$ws.Cells.Item(52, 7).Value = 8
$ws.Cells.Item(52, 8).Value = 533
$ws.Cells.Item(53, 2).Value = 21006
$ws.Cells.Item(53, 3).Value = 418
$ws.Cells.Item(53, 4).Value = 10144
$ws.Cells.Item(53, 5).Value = 10490
$ws.Cells.Item(53, 7).Value = 12
$ws.Cells.Item(53, 8).Value = 372
$ws.Cells.Item(55, 2).Value = 17968
$ws.Cells.Item(55, 3).Value = 52
$ws.Cells.Item(55, 4).Value = 16212
$ws.Cells.Item(55, 5).Value = 801
$ws.Cells.Item(55, 7).Value = 2
$ws.Cells.Item(55, 8).Value = 955
$ws.Cells.Item(67, 2).Value = 10650
$ws.Cells.Item(67, 3).Value = 127
$ws.Cells.Item(67, 5).Value = 2756
$ws.Cells.Item(70, 2).Value = 8796
$ws.Cells.Item(70, 3).Value = 98
$ws.Cells.Item(70, 4).Value = 3599
$ws.Cells.Item(70, 5).Value = 4656
$ws.Cells.Item(70, 7).Value = 8
$ws.Cells.Item(70, 8).Value = 541
$ws.Cells.Item(71, 2).Value = 8772
$ws.Cells.Item(71, 3).Value = 21
$ws.Cells.Item(71, 4).Value = 8138
$ws.Cells.Item(71, 5).Value = 386
$ws.Cells.Item(71, 8).Value = 248
$ws.Cells.Item(73, 2).Value = 7904
$ws.Cells.Item(73, 3).Value = 227
$ws.Cells.Item(73, 5).Value = 4353
$ws.Cells.Item(73, 7).Value = 2
$ws.Cells.Item(73, 8).Value = 58
$ws.Cells.Item(85, 2).Value = 4849
$ws.Cells.Item(85, 3).Value = 110
$ws.Cells.Item(85, 4).Value = 2107
$ws.Cells.Item(85, 5).Value = 2703
$ws.Cells.Item(85, 8).Value = 39
$ws.Cells.Item(86, 2).Value = 4848
$ws.Cells.Item(86, 3).Value = 185
$ws.Cells.Item(86, 4).Value = 1412
$ws.Cells.Item(86, 5).Value = 3361
$ws.Cells.Item(86, 8).Value = 75
$ws.Cells.Item(89, 2).Value = 4114
$ws.Cells.Item(89, 3).Value = 130
$ws.Cells.Item(89, 4).Value = 2217
$ws.Cells.Item(89, 5).Value = 1689
$ws.Cells.Item(89, 8).Value = 208
$ws.Cells.Item(90, 2).Value = 4107
$ws.Cells.Item(90, 3).Value = 5
$ws.Cells.Item(90, 4).Value = 2600
$ws.Cells.Item(90, 5).Value = 934
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 573
$ws.Cells.Item(91, 2).Value = 4048
$ws.Cells.Item(91, 4).Value = 1327
$ws.Cells.Item(91, 5).Value = 2686
$ws.Cells.Item(91, 8).Value = 35
$ws.Cells.Item(97, 2).Value = 3051
$ws.Cells.Item(97, 3).Value = 88
$ws.Cells.Item(97, 4).Value = 522
$ws.Cells.Item(97, 5).Value = 2492
$ws.Cells.Item(97, 7).Value = 7
$ws.Cells.Item(97, 8).Value = 37
$ws.Cells.Item(167, 4).Value = 154
$ws.Cells.Item(167, 5).Value = 40
$ws.Cells.Item(193, 4).Value = 27
$ws.Cells.Item(193, 5).Value = 2
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(212, 4).Value = 11
$ws.Cells.Item(212, 8).Value = 0

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 00:58"
